$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "54.011.16"
$ws.Range("E2").Value = "  +1.10%  "
$ws.Range("D3").Value = "2.241.24"
$ws.Range("E3").Value = "  +1.67%  "
$ws.Range("D4").Value = "'1.01"
$ws.Range("E4").Value = "  +0.47%  "
$ws.Range("D5").Value = "'495.40"
$ws.Range("E5").Value = "  +2.85%  "
$ws.Range("D6").Value = "'127.53"
$ws.Range("E6").Value = "  +2.07%  "
$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "  -0.33%  "
$ws.Range("D8").Value = "'0.527"
$ws.Range("E8").Value = "  +1.97%  "
$ws.Range("D9").Value = "2.279.64"
$ws.Range("E9").Value = "  +3.31%  "
$ws.Range("E10").Value = "  +4.47%  "
$ws.Range("E11").Value = "  +2.15%  "
$ws.Range("D12").Value = "'0.325"
$ws.Range("E12").Value = "  +3.96%  "
$ws.Range("D13").Value = "'4.63"
$ws.Range("E13").Value = "  -0.67%  "
$ws.Range("D14").Value = "2.652.56"
$ws.Range("E14").Value = "  +1.94%  "
$ws.Range("D15").Value = "'21.73"
$ws.Range("E15").Value = "  +3.86%  "
$ws.Range("D16").Value = "54.078.43"
$ws.Range("E16").Value = "  +1.35%  "
$ws.Range("E17").Value = "  +1.56%  "
$ws.Range("D18").Value = "2.300.31"
$ws.Range("E18").Value = "  +3.87%  "
$ws.Range("D19").Value = "'10.00"
$ws.Range("E19").Value = "  +5.36%  "
$ws.Range("D20").Value = "'4.09"
$ws.Range("E20").Value = "  +4.32%  "
$ws.Range("D21").Value = "'300.47"
$ws.Range("E21").Value = "  +1.40%  "
$ws.Range("D22").Value = "'6.42"
$ws.Range("E22").Value = "  +5.90%  "
$ws.Range("D23").Value = "'0.997"
$ws.Range("E23").Value = "  -0.40%  "
$ws.Range("E24").Value = "  -2.04%  "
$ws.Range("D25").Value = "'62.25"
$ws.Range("E25").Value = "  -0.97%  "
$ws.Range("D26").Value = "'1.02"
$ws.Range("E26").Value = "  +2.27%  "
$ws.Range("D27").Value = "2.396.06"
$ws.Range("E27").Value = "  +3.21%  "
$ws.Range("E28").Value = "  +3.08%  "
$ws.Range("E29").Value = "  +4.32%  "
$ws.Range("E30").Value = "  +1.70%  "
$ws.Range("D31").Value = "'168.39"
$ws.Range("E31").Value = "  +0.45%  "
$ws.Range("D32").Value = "'1.60"
$ws.Range("E32").Value = "  +2.41%  "
$ws.Range("D33").Value = "0.0₃0686"
$ws.Range("E33").Value = "  +2.67%  "
$ws.Range("D34").Value = "'5.86"
$ws.Range("E34").Value = "  +4.16%  "
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("D36").Value = "'0.991"
$ws.Range("E36").Value = "  -0.76%  "
$ws.Range("E37").Value = "  +2.34%  "
$ws.Range("D38").Value = "'17.73"
$ws.Range("E38").Value = "  +2.90%  "
$ws.Range("E39").Value = "  +10.08%  "
$ws.Range("E40").Value = "  +4.48%  "
$ws.Range("D41").Value = "'3.68"
$ws.Range("E41").Value = "  +4.34%  "
$ws.Range("E42").Value = "  -0.36%  "
$ws.Range("E43").Value = "  +3.59%  "
$ws.Range("E44").Value = "  +3.07%  "
$ws.Range("E45").Value = "  +3.84%  "
$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").Value = "'4.88"
$ws.Range("E46").Value = "  +5.97%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "'126.69"
$ws.Range("E47").Value = "  +4.01%  "
$ws.Range("D48").Value = "'0.0886"
$ws.Range("E48").Value = "  +1.42%  "
$ws.Range("E49").Value = "  +3.00%  "
$ws.Range("D50").Value = "'237.37"
$ws.Range("E50").Value = "  +3.68%  "
$ws.Range("E51").Value = "  +3.47%  "
